$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Restore C10 ("Integer max" value for rule R20) to 1
$ws.Range("C10").Value = 1
